$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "M1"
$ws.Cells.Item(2,2).Value = "Osm"
$ws.Cells.Item(2,3).Value = "Osmr"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 21.049038
$ws.Cells.Item(2,8).Value = 63.147114
$ws.Cells.Item(2,9).Value = 0.384846371905728
$ws.Cells.Item(2,10).Value = 0.384846371905728
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 50.640898
$ws.Cells.Item(2,14).Value = 151.922694
$ws.Cells.Item(2,15).Value = 0.2991745142287969
$ws.Cells.Item(2,16).Value = 0.3247737340214803
$ws.Cells.Item(2,17).Value = 1065.942186356124
$ws.Cells.Item(2,18).Value = 9593.479677205116
$ws.Cells.Item(2,19).Value = 0.1151362263676111
$ws.Cells.Item(2,20).Value = 0.1249879932284426

# Row 3
$ws.Cells.Item(3,1).Value = "M1"
$ws.Cells.Item(3,2).Value = "Osm"
$ws.Cells.Item(3,3).Value = "Osmr"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 21.049038
$ws.Cells.Item(3,8).Value = 63.147114
$ws.Cells.Item(3,9).Value = 0.384846371905728
$ws.Cells.Item(3,10).Value = 0.384846371905728
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 78.50335200000001
$ws.Cells.Item(3,14).Value = 235.510056
$ws.Cells.Item(3,15).Value = 0.4637793389827378
$ws.Cells.Item(3,16).Value = 0.5034631645403019
$ws.Cells.Item(3,17).Value = 1652.420039375376
$ws.Cells.Item(3,18).Value = 14871.78035437839
$ws.Cells.Item(3,19).Value = 0.1784837959723434
$ws.Cells.Item(3,20).Value = 0.1937559722615118

# Row 4
$ws.Cells.Item(4,1).Value = "M1"
$ws.Cells.Item(4,2).Value = "Osm"
$ws.Cells.Item(4,3).Value = "Osmr"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 21.049038
$ws.Cells.Item(4,8).Value = 63.147114
$ws.Cells.Item(4,9).Value = 0.384846371905728
$ws.Cells.Item(4,10).Value = 0.384846371905728
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.05322866666666667
$ws.Cells.Item(4,14).Value = 0.159686
$ws.Cells.Item(4,15).Value = 0.0003144624428470157
$ws.Cells.Item(4,16).Value = 0.0003413697922639135
$ws.Cells.Item(4,17).Value = 1.120412227356
$ws.Cells.Item(4,18).Value = 10.083710046204
$ws.Cells.Item(4,19).Value = 0.0001210197302302863
$ws.Cells.Item(4,20).Value = 0.0001313749260309792

# Row 5
$ws.Cells.Item(5,1).Value = "M1"
$ws.Cells.Item(5,2).Value = "Osm"
$ws.Cells.Item(5,3).Value = "Osmr"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 21.049038
$ws.Cells.Item(5,8).Value = 63.147114
$ws.Cells.Item(5,9).Value = 0.384846371905728
$ws.Cells.Item(5,10).Value = 0.384846371905728
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.04512133333333334
$ws.Cells.Item(5,14).Value = 0.135364
$ws.Cells.Item(5,15).Value = 0.0002665662244250807
$ws.Cells.Item(5,16).Value = 0.0002893752774821361
$ws.Cells.Item(5,17).Value = 0.9497606599440002
$ws.Cells.Item(5,18).Value = 8.547845939496002
$ws.Cells.Item(5,19).Value = 0.0001025870443426004
$ws.Cells.Item(5,20).Value = 0.0001113650256582134

# Row 6
$ws.Cells.Item(6,1).Value = "M1"
$ws.Cells.Item(6,2).Value = "Osm"
$ws.Cells.Item(6,3).Value = "Osmr"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 21.049038
$ws.Cells.Item(6,8).Value = 63.147114
$ws.Cells.Item(6,9).Value = 0.384846371905728
$ws.Cells.Item(6,10).Value = 0.384846371905728
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 40.0261565
$ws.Cells.Item(6,14).Value = 80.052313
$ws.Cells.Item(6,15).Value = 0.2364651181211933
$ws.Cells.Item(6,16).Value = 0.1711323563684718
$ws.Cells.Item(6,17).Value = 842.512089162447
$ws.Cells.Item(6,18).Value = 5055.072534974682
$ws.Cells.Item(6,19).Value = 0.09100274279120063
$ws.Cells.Item(6,20).Value = 0.06585966646408446

# Row 7
$ws.Cells.Item(7,1).Value = "M2"
$ws.Cells.Item(7,2).Value = "Osm"
$ws.Cells.Item(7,3).Value = "Osmr"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 33.645613
$ws.Cells.Item(7,8).Value = 100.936839
$ws.Cells.Item(7,9).Value = 0.615153628094272
$ws.Cells.Item(7,10).Value = 0.615153628094272
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 50.640898
$ws.Cells.Item(7,14).Value = 151.922694
$ws.Cells.Item(7,15).Value = 0.2991745142287969
$ws.Cells.Item(7,16).Value = 0.3247737340214803
$ws.Cells.Item(7,17).Value = 1703.844056080474
$ws.Cells.Item(7,18).Value = 15334.59650472427
$ws.Cells.Item(7,19).Value = 0.1840382878611858
$ws.Cells.Item(7,20).Value = 0.1997857407930377

# Row 8
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Osm"
$ws.Cells.Item(8,3).Value = "Osmr"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 33.645613
$ws.Cells.Item(8,8).Value = 100.936839
$ws.Cells.Item(8,9).Value = 0.615153628094272
$ws.Cells.Item(8,10).Value = 0.615153628094272
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 78.50335200000001
$ws.Cells.Item(8,14).Value = 235.510056
$ws.Cells.Item(8,15).Value = 0.4637793389827378
$ws.Cells.Item(8,16).Value = 0.5034631645403019
$ws.Cells.Item(8,17).Value = 2641.293400594776
$ws.Cells.Item(8,18).Value = 23771.64060535299
$ws.Cells.Item(8,19).Value = 0.2852955430103944
$ws.Cells.Item(8,20).Value = 0.3097071922787902

# Row 9
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Osm"
$ws.Cells.Item(9,3).Value = "Osmr"
$ws.Cells.Item(9,4).Value = "M1"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 33.645613
$ws.Cells.Item(9,8).Value = 100.936839
$ws.Cells.Item(9,9).Value = 0.615153628094272
$ws.Cells.Item(9,10).Value = 0.615153628094272
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.05322866666666667
$ws.Cells.Item(9,14).Value = 0.159686
$ws.Cells.Item(9,15).Value = 0.0003144624428470157
$ws.Cells.Item(9,16).Value = 0.0003413697922639135
$ws.Cells.Item(9,17).Value = 1.790911119172667
$ws.Cells.Item(9,18).Value = 16.118200072554
$ws.Cells.Item(9,19).Value = 0.0001934427126167294
$ws.Cells.Item(9,20).Value = 0.0002099948662329343

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Osm"
$ws.Cells.Item(10,3).Value = "Osmr"
$ws.Cells.Item(10,4).Value = "M2"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 33.645613
$ws.Cells.Item(10,8).Value = 100.936839
$ws.Cells.Item(10,9).Value = 0.615153628094272
$ws.Cells.Item(10,10).Value = 0.615153628094272
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.04512133333333334
$ws.Cells.Item(10,14).Value = 0.135364
$ws.Cells.Item(10,15).Value = 0.0002665662244250807
$ws.Cells.Item(10,16).Value = 0.0002893752774821361
$ws.Cells.Item(10,17).Value = 1.518134919377334
$ws.Cells.Item(10,18).Value = 13.663214274396
$ws.Cells.Item(10,19).Value = 0.0001639791800824803
$ws.Cells.Item(10,20).Value = 0.0001780102518239227

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Osm"
$ws.Cells.Item(11,3).Value = "Osmr"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 33.645613
$ws.Cells.Item(11,8).Value = 100.936839
$ws.Cells.Item(11,9).Value = 0.615153628094272
$ws.Cells.Item(11,10).Value = 0.615153628094272
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 40.0261565
$ws.Cells.Item(11,14).Value = 80.052313
$ws.Cells.Item(11,15).Value = 0.2364651181211933
$ws.Cells.Item(11,16).Value = 0.1711323563684718
$ws.Cells.Item(11,17).Value = 1346.704571476435
$ws.Cells.Item(11,18).Value = 8080.227428858608
$ws.Cells.Item(11,19).Value = 0.1454623753299926
$ws.Cells.Item(11,20).Value = 0.1052726899043873
